# Updated Product Backlog with Priorities
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sprint 1")

# --- Task Description (column D) content updates ---
$ws.Range("D10").Value = "Displays pictures taken by the camera in History page."
$ws.Range("D11").Value = "Profile of admins should contain: username, email, password, identifier of door, and a list of admins who have access to the same door."
$ws.Range("D12").Value = "The presenter should be able to create a profile if they don't have one and be able to login once they have one by entering an door identifier, email address, password, and username."

# --- Comments (column H) - append Priority lines to "Assigned to:" notes ---
$ws.Range("H7").Value = "Assigned to: Irfan Ahmed`nPriority 1"
$ws.Range("H8").Value = "Assigned to: Irfan Ahmed`nPriority 5"
$ws.Range("H9").Value = "Assigned to: Cong-Vinh Vu`nPriority 3"
$ws.Range("H10").Value = "Assigned to: Alec Kurkdjian`nPriority 3"
$ws.Range("H11").Value = "Assigned to: Alec Kurkdjian`nPriority 1"
$ws.Range("H12").Value = "Assigned to: Cong-Vinh Vu`nPriority 2"
$ws.Range("H13").Value = "Assigned to: Cong-Vinh Vu`nPriority 2"
$ws.Range("H14").Value = "Assigned to: Alec Kurkdjian and Cong-Vinh Vu`nPriority 1"
$ws.Range("H15").Value = "Assigned to: Alec Kurkdjian and Cong-Vinh Vu`nPriority 2"
$ws.Range("H16").Value = "Assigned to: Alec Kurkdjian and Cong-Vinh Vu`nPriority 1"
$ws.Range("H17").Value = "Assigned to: Shadi Makdissi and Ogo-Oluwa Jesutomi Olasubulumi`n"
$ws.Range("H18").Value = "Assigned to: Shadi Makdissi and Ogo-Oluwa Jesutomi Olasubulumi`n"

# Newly multi-line "Assigned to:" comments need word-wrap turned on to match
# the rest of the column (the others already wrap).
$ws.Range("H7").WrapText = $true
$ws.Range("H8").WrapText = $true
$ws.Range("H9").WrapText = $true

# --- Row heights follow the new wrapped line counts ---
$ws.Rows.Item(10).RowHeight = 28.8
$ws.Rows.Item(11).RowHeight = 43.2
$ws.Rows.Item(12).RowHeight = 57.6
$ws.Rows.Item(14).RowHeight = 28.8
$ws.Rows.Item(15).RowHeight = 28.8
$ws.Rows.Item(17).RowHeight = 43.2
$ws.Rows.Item(18).RowHeight = 43.2

# --- Restore the selection left after editing ---
$ws.Range("D12").Select()
